$d = $word.ActiveDocument

# --- Locate a paragraph that already uses the numId=8 list (List Paragraph
#     style + numbering) so we can reuse/continue that same list instance.
$srcListPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Implementación de un constructor de copias*") {
        $srcListPara = $p
        break
    }
}
$listTemplate = $srcListPara.Range.ListFormat.ListTemplate

# --- Change 1: table cell that only contained "... " becomes a
#     single-item bulleted/numbered ("Prrafodelista" / numId 8) paragraph
#     reading "Escribir código".
$found = $d.Content.Find.Execute("... ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Escribir código", 2)

$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Escribir código*" -and $p.Range.Text -notlike "*constructor*") {
        $targetPara = $p
    }
}
$targetPara.Style = "List Paragraph"
$targetPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 1, $false, 1)

# --- Change 2: append two new list items after "Creación de verificador..."
$verifPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Creación de verificador así como de seguro para inversa de matrices cuadradas*") {
        $verifPara = $p
    }
}

$verifPara.Range.InsertParagraphAfter()
$sumaPara = $verifPara.Next()
$sumaPara.Range.Text = "Creación del método que suma dos vectores"
$sumaPara.Style = "List Paragraph"
$sumaPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 1, $false, 1)

$sumaPara.Range.InsertParagraphAfter()
$restaPara = $sumaPara.Next()
$restaPara.Range.Text = "Creación del método que resta dos vectores"
$restaPara.Style = "List Paragraph"
$restaPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 1, $false, 1)

Write-Host "Edits applied."
